$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 13913
$ws.Cells.Item(2, 5).Value = 44
$ws.Cells.Item(2, 6).Value = 44
$ws.Cells.Item(2, 7).Value = 190
$ws.Cells.Item(2, 8).Value = 133
$ws.Cells.Item(2, 9).Value = 104
$ws.Cells.Item(2, 10).Value = 28
$ws.Cells.Item(2, 11).Value = 13466
$ws.Cells.Item(2, 12).Value = 7128
$ws.Cells.Item(2, 13).Value = 6338
$ws.Cells.Item(2, 14).Value = 5809
$ws.Cells.Item(2, 15).Value = 529
$ws.Cells.Item(2, 16).Value = 300
$ws.Cells.Item(2, 17).Value = 368
$ws.Cells.Item(2, 18).Value = -256
$ws.Cells.Item(2, 19).Value = -146
$ws.Cells.Item(2, 20).Value = 271
$ws.Cells.Item(2, 21).Value = 97
$ws.Cells.Item(2, 22).Value = 2039
$ws.Cells.Item(2, 23).Value = 0.31
$ws.Cells.Item(2, 24).Value = 0.95
$ws.Cells.Item(2, 25).Value = 1.79
$ws.Cells.Item(2, 26).Value = 0.98
$ws.Cells.Item(2, 27).Value = 112.47
$ws.Cells.Item(2, 28).Value = 1873.17
$ws.Cells.Item(2, 29).Value = 1737
$ws.Cells.Item(2, 30).Value = 20.79
$ws.Cells.Item(2, 31).Value = 115028
$ws.Cells.Item(2, 32).Value = 0.31
$ws.Cells.Item(2, 33).Value = 1250
$ws.Cells.Item(2, 34).Value = 3.46
$ws.Cells.Item(2, 35).Value = 60.59
$ws.Cells.Item(2, 36).Value = 6000000

# Row 3
$ws.Cells.Item(3, 4).Value = 12007
$ws.Cells.Item(3, 5).Value = 174
$ws.Cells.Item(3, 6).Value = 174
$ws.Cells.Item(3, 7).Value = 229
$ws.Cells.Item(3, 8).Value = 194
$ws.Cells.Item(3, 9).Value = 158
$ws.Cells.Item(3, 10).Value = 36
$ws.Cells.Item(3, 11).Value = 12646
$ws.Cells.Item(3, 12).Value = 6173
$ws.Cells.Item(3, 13).Value = 6473
$ws.Cells.Item(3, 14).Value = 5909
$ws.Cells.Item(3, 15).Value = 565
$ws.Cells.Item(3, 16).Value = 300
$ws.Cells.Item(3, 17).Value = 546
$ws.Cells.Item(3, 18).Value = -1034
$ws.Cells.Item(3, 19).Value = -119
$ws.Cells.Item(3, 20).Value = 692
$ws.Cells.Item(3, 21).Value = -146
$ws.Cells.Item(3, 22).Value = 1995
$ws.Cells.Item(3, 23).Value = 1.45
$ws.Cells.Item(3, 24).Value = 1.62
$ws.Cells.Item(3, 25).Value = 2.7
$ws.Cells.Item(3, 26).Value = 1.49
$ws.Cells.Item(3, 27).Value = 95.34999999999999
$ws.Cells.Item(3, 28).Value = 1894.66
$ws.Cells.Item(3, 29).Value = 2633
$ws.Cells.Item(3, 30).Value = 13.95
$ws.Cells.Item(3, 31).Value = 116997
$ws.Cells.Item(3, 32).Value = 0.31
$ws.Cells.Item(3, 33).Value = 1250
$ws.Cells.Item(3, 34).Value = 3.4
$ws.Cells.Item(3, 35).Value = 39.95
$ws.Cells.Item(3, 36).Value = 6000000

# Row 4
$ws.Cells.Item(4, 4).Value = 10513
$ws.Cells.Item(4, 5).Value = 211
$ws.Cells.Item(4, 6).Value = 211
$ws.Cells.Item(4, 7).Value = 562
$ws.Cells.Item(4, 8).Value = 492
$ws.Cells.Item(4, 9).Value = 351
$ws.Cells.Item(4, 10).Value = 141
$ws.Cells.Item(4, 11).Value = 12861
$ws.Cells.Item(4, 12).Value = 6114
$ws.Cells.Item(4, 13).Value = 6747
$ws.Cells.Item(4, 14).Value = 6153
$ws.Cells.Item(4, 15).Value = 594
$ws.Cells.Item(4, 16).Value = 300
$ws.Cells.Item(4, 17).Value = 355
$ws.Cells.Item(4, 18).Value = -941
$ws.Cells.Item(4, 19).Value = -100
$ws.Cells.Item(4, 20).Value = 396
$ws.Cells.Item(4, 21).Value = -41
$ws.Cells.Item(4, 22).Value = 2056
$ws.Cells.Item(4, 23).Value = 2.01
$ws.Cells.Item(4, 24).Value = 4.68
$ws.Cells.Item(4, 25).Value = 5.81
$ws.Cells.Item(4, 26).Value = 3.86
$ws.Cells.Item(4, 27).Value = 90.62
$ws.Cells.Item(4, 28).Value = 1996.46
$ws.Cells.Item(4, 29).Value = 5842
$ws.Cells.Item(4, 30).Value = 6.35
$ws.Cells.Item(4, 31).Value = 121833
$ws.Cells.Item(4, 32).Value = 0.3
$ws.Cells.Item(4, 33).Value = 1500
$ws.Cells.Item(4, 34).Value = 4.04
$ws.Cells.Item(4, 35).Value = 21.61
$ws.Cells.Item(4, 36).Value = 6000000

# Row 5
$ws.Cells.Item(5, 4).Value = 11152
$ws.Cells.Item(5, 5).Value = 244
$ws.Cells.Item(5, 6).Value = 244
$ws.Cells.Item(5, 7).Value = 422
$ws.Cells.Item(5, 8).Value = 271
$ws.Cells.Item(5, 9).Value = 226
$ws.Cells.Item(5, 10).Value = 44
$ws.Cells.Item(5, 11).Value = 14350
$ws.Cells.Item(5, 12).Value = 7376
$ws.Cells.Item(5, 13).Value = 6974
$ws.Cells.Item(5, 14).Value = 6392
$ws.Cells.Item(5, 15).Value = 582
$ws.Cells.Item(5, 16).Value = 300
$ws.Cells.Item(5, 17).Value = 848
$ws.Cells.Item(5, 18).Value = -1063
$ws.Cells.Item(5, 19).Value = 632
$ws.Cells.Item(5, 20).Value = 143
$ws.Cells.Item(5, 21).Value = 706
$ws.Cells.Item(5, 22).Value = 2796
$ws.Cells.Item(5, 23).Value = 2.19
$ws.Cells.Item(5, 24).Value = 2.43
$ws.Cells.Item(5, 25).Value = 3.61
$ws.Cells.Item(5, 26).Value = 1.99
$ws.Cells.Item(5, 27).Value = 105.75
$ws.Cells.Item(5, 28).Value = 2054.72
$ws.Cells.Item(5, 29).Value = 3772
$ws.Cells.Item(5, 30).Value = 10.05
$ws.Cells.Item(5, 31).Value = 126563
$ws.Cells.Item(5, 32).Value = 0.3
$ws.Cells.Item(5, 33).Value = 1500
$ws.Cells.Item(5, 34).Value = 3.96
$ws.Cells.Item(5, 35).Value = 33.47
$ws.Cells.Item(5, 36).Value = 6000000

# Row 6
$ws.Cells.Item(6, 4).Value = 10954
$ws.Cells.Item(6, 5).Value = 252
$ws.Cells.Item(6, 6).Value = 252
$ws.Cells.Item(6, 7).Value = 85
$ws.Cells.Item(6, 8).Value = 107
$ws.Cells.Item(6, 9).Value = 70
$ws.Cells.Item(6, 11).Value = 15206
$ws.Cells.Item(6, 12).Value = 8736
$ws.Cells.Item(6, 13).Value = 6470
$ws.Cells.Item(6, 14).Value = 5862
$ws.Cells.Item(6, 16).Value = 300
$ws.Cells.Item(6, 17).Value = -212
$ws.Cells.Item(6, 18).Value = -280
$ws.Cells.Item(6, 19).Value = 775
$ws.Cells.Item(6, 20).Value = 234
$ws.Cells.Item(6, 21).Value = -446
$ws.Cells.Item(6, 22).Value = 4010
$ws.Cells.Item(6, 23).Value = 2.3
$ws.Cells.Item(6, 24).Value = 0.98
$ws.Cells.Item(6, 25).Value = 1.15
$ws.Cells.Item(6, 26).Value = 0.72
$ws.Cells.Item(6, 27).Value = 135.01
$ws.Cells.Item(6, 28).Value = 2020.72
$ws.Cells.Item(6, 29).Value = 1173
$ws.Cells.Item(6, 30).Value = 34.05
$ws.Cells.Item(6, 31).Value = 137418
$ws.Cells.Item(6, 32).Value = 0.29
$ws.Cells.Item(6, 33).Value = 1500
$ws.Cells.Item(6, 34).Value = 3.75
$ws.Cells.Item(6, 35).Value = 90.90000000000001
$ws.Cells.Item(6, 36).Value = 6000000

# Row 7 - clear cells
$ws.Cells.Item(7, 4).ClearContents()
$ws.Cells.Item(7, 5).ClearContents()
$ws.Cells.Item(7, 7).ClearContents()
$ws.Cells.Item(7, 8).ClearContents()
$ws.Cells.Item(7, 9).ClearContents()
$ws.Cells.Item(7, 11).ClearContents()
$ws.Cells.Item(7, 12).ClearContents()
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(7, 14).ClearContents()
$ws.Cells.Item(7, 16).ClearContents()
$ws.Cells.Item(7, 17).ClearContents()
$ws.Cells.Item(7, 18).ClearContents()
$ws.Cells.Item(7, 19).ClearContents()
$ws.Cells.Item(7, 20).ClearContents()
$ws.Cells.Item(7, 21).ClearContents()
$ws.Cells.Item(7, 23).ClearContents()
$ws.Cells.Item(7, 24).ClearContents()
$ws.Cells.Item(7, 25).ClearContents()
$ws.Cells.Item(7, 26).ClearContents()
$ws.Cells.Item(7, 27).ClearContents()
$ws.Cells.Item(7, 29).ClearContents()
$ws.Cells.Item(7, 30).ClearContents()
$ws.Cells.Item(7, 31).ClearContents()
$ws.Cells.Item(7, 32).ClearContents()
$ws.Cells.Item(7, 33).ClearContents()
$ws.Cells.Item(7, 34).ClearContents()
$ws.Cells.Item(7, 35).ClearContents()

# Row 8 - clear cells
$ws.Cells.Item(8, 4).ClearContents()
$ws.Cells.Item(8, 5).ClearContents()
$ws.Cells.Item(8, 7).ClearContents()
$ws.Cells.Item(8, 8).ClearContents()
$ws.Cells.Item(8, 9).ClearContents()
$ws.Cells.Item(8, 11).ClearContents()
$ws.Cells.Item(8, 12).ClearContents()
$ws.Cells.Item(8, 13).ClearContents()
$ws.Cells.Item(8, 14).ClearContents()
$ws.Cells.Item(8, 16).ClearContents()
$ws.Cells.Item(8, 17).ClearContents()
$ws.Cells.Item(8, 18).ClearContents()
$ws.Cells.Item(8, 19).ClearContents()
$ws.Cells.Item(8, 20).ClearContents()
$ws.Cells.Item(8, 21).ClearContents()
$ws.Cells.Item(8, 23).ClearContents()
$ws.Cells.Item(8, 24).ClearContents()
$ws.Cells.Item(8, 25).ClearContents()
$ws.Cells.Item(8, 26).ClearContents()
$ws.Cells.Item(8, 27).ClearContents()
$ws.Cells.Item(8, 29).ClearContents()
$ws.Cells.Item(8, 30).ClearContents()
$ws.Cells.Item(8, 31).ClearContents()
$ws.Cells.Item(8, 32).ClearContents()
$ws.Cells.Item(8, 33).ClearContents()
$ws.Cells.Item(8, 34).ClearContents()
$ws.Cells.Item(8, 35).ClearContents()

# Row 9 - clear cells
$ws.Cells.Item(9, 4).ClearContents()
$ws.Cells.Item(9, 5).ClearContents()
$ws.Cells.Item(9, 7).ClearContents()
$ws.Cells.Item(9, 8).ClearContents()
$ws.Cells.Item(9, 9).ClearContents()
$ws.Cells.Item(9, 11).ClearContents()
$ws.Cells.Item(9, 12).ClearContents()
$ws.Cells.Item(9, 13).ClearContents()
$ws.Cells.Item(9, 14).ClearContents()
$ws.Cells.Item(9, 16).ClearContents()
$ws.Cells.Item(9, 17).ClearContents()
$ws.Cells.Item(9, 18).ClearContents()
$ws.Cells.Item(9, 19).ClearContents()
$ws.Cells.Item(9, 20).ClearContents()
$ws.Cells.Item(9, 21).ClearContents()
$ws.Cells.Item(9, 23).ClearContents()
$ws.Cells.Item(9, 24).ClearContents()
$ws.Cells.Item(9, 25).ClearContents()
$ws.Cells.Item(9, 26).ClearContents()
$ws.Cells.Item(9, 27).ClearContents()
$ws.Cells.Item(9, 29).ClearContents()
$ws.Cells.Item(9, 30).ClearContents()
$ws.Cells.Item(9, 31).ClearContents()
$ws.Cells.Item(9, 32).ClearContents()
$ws.Cells.Item(9, 33).ClearContents()
$ws.Cells.Item(9, 34).ClearContents()
$ws.Cells.Item(9, 35).ClearContents()
